$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (pushes old rows 9-44 down to 10-45)
$ws.Rows.Item(9).Insert()

# New row 9: "START" marker mirroring the header row's formatting
$ws.Range("A8:J8").Copy()
$ws.Range("A9:J9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(9, 1).Value = "START"

# The spacer row that shifted down to row 21 gets an "END" marker
$ws.Cells.Item(21, 1).Value = "END"

$ws.Range("A21").Select()
